$wb = $excel.ActiveWorkbook

# --- Add the new worksheet "Best selling cars in 2022 ww" at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Best selling cars in 2022 ww"

# --- Populate Model column first (rows 3-12), top to bottom ---
$newSheet.Range("B3").Value = "Toyota Corolla"
$newSheet.Range("B4").Value = "Toyota RAV4"
$newSheet.Range("B5").Value = "Ford F-Series"
$newSheet.Range("B6").Value = "Tesla Model Y"
$newSheet.Range("B7").Value = "Toyota Camry"
$newSheet.Range("B8").Value = "Honda CR-V"
$newSheet.Range("B9").Value = "Tesla Model 3"
$newSheet.Range("B10").Value = "Chevy Silverado"
$newSheet.Range("B11").Value = "Toyota Hilux"
$newSheet.Range("B12").Value = "Hyundai Tucson"

# --- Populate header row (Model, Brand, Units Sold) ---
$newSheet.Range("B2").Value = "Model"
$newSheet.Range("A2").Value = "Brand"
$newSheet.Range("C2").Value = "Units Sold"

# --- Populate Brand column (rows 3-12), visiting distinct brands in first-use order ---
$newSheet.Range("A3").Value = "Toyota"
$newSheet.Range("A5").Value = "Ford"
$newSheet.Range("A6").Value = "Tesla"
$newSheet.Range("A8").Value = "Honda"
$newSheet.Range("A12").Value = "Hyundai"
$newSheet.Range("A10").Value = "Chervolet"
$newSheet.Range("A4").Value = "Toyota"
$newSheet.Range("A7").Value = "Toyota"
$newSheet.Range("A11").Value = "Toyota"
$newSheet.Range("A9").Value = "Tesla"

# --- Populate Units Sold column (numbers, order irrelevant) ---
$newSheet.Range("C3").Value = 1120000
$newSheet.Range("C4").Value = 870000
$newSheet.Range("C5").Value = 787000
$newSheet.Range("C6").Value = 786000
$newSheet.Range("C7").Value = 675000
$newSheet.Range("C8").Value = 601000
$newSheet.Range("C9").Value = 596000
$newSheet.Range("C10").Value = 592000
$newSheet.Range("C11").Value = 564000
$newSheet.Range("C12").Value = 564000

# --- Column widths for new sheet (best-fit approximation) ---
$newSheet.Columns.Item(1).ColumnWidth = 8.88
$newSheet.Columns.Item(2).ColumnWidth = 14.45

# --- Selection / view for new sheet ---
$newSheet.Range("H17").Select()

# --- Selection on "Copied" sheet ---
$copiedSheet = $wb.Worksheets.Item("Copied")
$copiedSheet.Range("J45").Select()

# --- Reset (auto-fit) row heights on "Country Sales" that previously had explicit heights ---
$countrySheet = $wb.Worksheets.Item("Country Sales")
$countrySheet.Rows.Item(2).EntireRow.AutoFit()
$countrySheet.Rows.Item(3).EntireRow.AutoFit()
$countrySheet.Rows.Item(9).EntireRow.AutoFit()
$countrySheet.Rows.Item(17).EntireRow.AutoFit()
$countrySheet.Rows.Item(19).EntireRow.AutoFit()
$countrySheet.Rows.Item(24).EntireRow.AutoFit()

# --- topLeftCell scroll position on "Toyota manufactoring plants" ---
$plantsSheet = $wb.Worksheets.Item("Toyota manufactoring plants")
$plantsSheet.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1

# --- Activate the new sheet last so it becomes the selected tab ---
$newSheet.Activate()

# --- Adjust workbook window size/position ---
$excel.ActiveWindow.WindowState = -4143
$excel.Width = 1200
$excel.Height = 660
